$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 40928.16
$ws.Range("I28").Value = 60087.06
$ws.Range("J28").Value = 215.5
$ws.Range("K28").Value = 60087.06
$ws.Range("L28").Value = 215.5
$ws.Range("M28").Value = -59602.06
$ws.Range("N28").Value = -1185.5
# Row 116
$ws.Range("H116").Value = 10909.4375
$ws.Range("I116").Value = 3718.375
$ws.Range("K116").Value = 3718.375
$ws.Range("M116").Value = -276.375
# Row 125
$ws.Range("H125").Value = 9263352
$ws.Range("I125").Value = 2022.25
$ws.Range("J125").Value = 13894017
$ws.Range("K125").Value = 18200.25
$ws.Range("L125").Value = 125046153
$ws.Range("M125").Value = -15740.25
$ws.Range("N125").Value = -125051073
# Row 132
$ws.Range("H132").Value = 2722.8948
$ws.Range("I132").Value = 3003.5
$ws.Range("J132").Value = 1226.3334
$ws.Range("K132").Value = 9010.5
$ws.Range("L132").Value = 3679.0002
$ws.Range("M132").Value = -6480.5
$ws.Range("N132").Value = -8739.0002
# Row 137
$ws.Range("H137").Value = 3271.6667
$ws.Range("I137").Value = 3007.4
$ws.Range("K137").Value = 9022.200000000001
$ws.Range("M137").Value = -6472.200000000001
# Row 138
$ws.Range("H138").Value = 4568.1665
$ws.Range("I138").Value = 2086.3333
$ws.Range("J138").Value = 5064.533
$ws.Range("K138").Value = 6258.999899999999
$ws.Range("L138").Value = 15193.599
$ws.Range("M138").Value = -1118.999899999999
$ws.Range("N138").Value = -25473.599

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3943.0706
$ws.Range("I32").Value = 3943.0706
$ws.Range("K32").Value = 3943.0706
$ws.Range("M32").Value = -3656.0706
# Row 45
$ws.Range("H45").Value = 1967
$ws.Range("I45").Value = 1961.5
$ws.Range("K45").Value = 1961.5
$ws.Range("M45").Value = -1584.5
# Row 61
$ws.Range("H61").Value = 2540.7368
$ws.Range("I61").Value = 2293
$ws.Range("K61").Value = 2293
$ws.Range("M61").Value = -2081
# Row 97
$ws.Range("H97").Value = 1500.0834
$ws.Range("I97").Value = 1377.9
$ws.Range("J97").Value = 2111
$ws.Range("K97").Value = 1377.9
$ws.Range("L97").Value = 2111
$ws.Range("M97").Value = -881.9000000000001
$ws.Range("N97").Value = -3103
# Row 110
$ws.Range("H110").Value = 219137.12
$ws.Range("I110").Value = 239906.72
$ws.Range("K110").Value = 239906.72
$ws.Range("M110").Value = -237861.72
# Row 122
$ws.Range("H122").Value = 5836.3076
$ws.Range("I122").Value = 6728
$ws.Range("J122").Value = 5279
$ws.Range("K122").Value = 20184
$ws.Range("L122").Value = 15837
$ws.Range("M122").Value = -17734
$ws.Range("N122").Value = -20737
# Row 132
$ws.Range("H132").Value = 2520.1396
$ws.Range("I132").Value = 2575.0527
$ws.Range("K132").Value = 7725.158100000001
$ws.Range("M132").Value = -5195.158100000001
# Row 136
$ws.Range("H136").Value = 2540.7368
$ws.Range("I136").Value = 2293
$ws.Range("K136").Value = 6879
$ws.Range("M136").Value = -4329

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1378.625
$ws.Range("I99").Value = 1171.6666
$ws.Range("J99").Value = 1999.5
$ws.Range("K99").Value = 1171.6666
$ws.Range("L99").Value = 1999.5
$ws.Range("M99").Value = 326.3334
$ws.Range("N99").Value = -4995.5
# Row 107
$ws.Range("H107").Value = 1625
$ws.Range("I107").Value = 978.0833
$ws.Range("K107").Value = 978.0833
$ws.Range("M107").Value = 941.9167
# Row 134
$ws.Range("H134").Value = 41977.383
$ws.Range("I134").Value = 3244
$ws.Range("K134").Value = 9732
$ws.Range("M134").Value = -7197

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1779.1765
$ws.Range("I16").Value = 1821.1428
$ws.Range("J16").Value = 1583.3334
$ws.Range("K16").Value = 1821.1428
$ws.Range("L16").Value = 1583.3334
$ws.Range("M16").Value = -1534.1428
$ws.Range("N16").Value = -2157.3334
# Row 20
$ws.Range("H20").Value = 66836.14
$ws.Range("J20").Value = 66836.14
$ws.Range("L20").Value = 66836.14
$ws.Range("N20").Value = -67308.14
# Row 30
$ws.Range("H30").Value = 66836.14
$ws.Range("J30").Value = 66836.14
$ws.Range("L30").Value = 66836.14
$ws.Range("N30").Value = -67018.14
# Row 94
$ws.Range("H94").Value = 3478
$ws.Range("I94").Value = 3012
$ws.Range("K94").Value = 3012
$ws.Range("M94").Value = -2561
# Row 107
$ws.Range("H107").Value = 557
$ws.Range("I107").Value = 343.14285
$ws.Range("K107").Value = 343.14285
$ws.Range("M107").Value = 1576.85715
# Row 113
$ws.Range("H113").Value = 1779.1765
$ws.Range("I113").Value = 1821.1428
$ws.Range("J113").Value = 1583.3334
$ws.Range("K113").Value = 1821.1428
$ws.Range("L113").Value = 1583.3334
$ws.Range("M113").Value = 348.8571999999999
$ws.Range("N113").Value = -5923.3334
# Row 128
$ws.Range("H128").Value = 66836.14
$ws.Range("J128").Value = 66836.14
$ws.Range("L128").Value = 66836.14
$ws.Range("N128").Value = -76796.14
# Row 132
$ws.Range("H132").Value = 1559.875
$ws.Range("I132").Value = 1496.6666
$ws.Range("J132").Value = 1749.5
$ws.Range("K132").Value = 4489.9998
$ws.Range("L132").Value = 5248.5
$ws.Range("M132").Value = -1959.9998
$ws.Range("N132").Value = -10308.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 360969.6
$ws.Range("I5").Value = 62367
$ws.Range("J5").Value = 589312.75
$ws.Range("K5").Value = 187101
$ws.Range("L5").Value = 1767938.25
$ws.Range("M5").Value = -186989
$ws.Range("N5").Value = -1768162.25
# Row 12
$ws.Range("H12").Value = 121.63158
$ws.Range("J12").Value = 135.82353
$ws.Range("L12").Value = 407.47059
$ws.Range("N12").Value = -753.47059
# Row 129
$ws.Range("H129").Value = 27862346
$ws.Range("I129").Value = 41667350
$ws.Range("J129").Value = 252333.25
$ws.Range("K129").Value = 125002050
$ws.Range("L129").Value = 756999.75
$ws.Range("M129").Value = -124997050
$ws.Range("N129").Value = -766999.75
# Row 135
$ws.Range("H135").Value = 360969.6
$ws.Range("I135").Value = 62367
$ws.Range("J135").Value = 589312.75
$ws.Range("K135").Value = 561303
$ws.Range("L135").Value = 5303814.75
$ws.Range("M135").Value = -558768
$ws.Range("N135").Value = -5308884.75

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 279
$ws.Range("I2").Value = 285.6842
$ws.Range("J2").Value = 236.66667
$ws.Range("K2").Value = 285.6842
$ws.Range("L2").Value = 236.66667
$ws.Range("M2").Value = -172.6842
$ws.Range("N2").Value = -462.66667
# Row 11
$ws.Range("H11").Value = 6683333.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6683333.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 6683333.5
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -6683611.5
# Row 113
$ws.Range("H113").Value = 508278.1
$ws.Range("J113").Value = 10224.2
$ws.Range("L113").Value = 10224.2
$ws.Range("N113").Value = -14564.2
# Row 122
$ws.Range("H122").Value = 3371.375
$ws.Range("I122").Value = 1678
$ws.Range("J122").Value = 5548.5713
$ws.Range("K122").Value = 5034
$ws.Range("L122").Value = 16645.7139
$ws.Range("M122").Value = -2584
$ws.Range("N122").Value = -21545.7139
# Row 132
$ws.Range("H132").Value = 47648.824
$ws.Range("I132").Value = 5307.4
$ws.Range("J132").Value = 127039
$ws.Range("K132").Value = 15922.2
$ws.Range("L132").Value = 381117
$ws.Range("M132").Value = -13392.2
$ws.Range("N132").Value = -386177

$ws = $wb.Worksheets.Item("LTW")
# Row 23
$ws.Range("H23").Value = 459799.6
$ws.Range("I23").Value = 459799.6
$ws.Range("K23").Value = 459799.6
$ws.Range("M23").Value = -459569.6
# Row 55
$ws.Range("H55").Value = 723.8387
$ws.Range("I55").Value = 220.19048
$ws.Range("K55").Value = 220.19048
$ws.Range("M55").Value = -47.19048000000001
# Row 68
$ws.Range("H68").Value = 253688.75
$ws.Range("I68").Value = 3002
$ws.Range("K68").Value = 3002
$ws.Range("M68").Value = -2253
# Row 71
$ws.Range("H71").Value = 253688.75
$ws.Range("I71").Value = 3002
$ws.Range("K71").Value = 15010
$ws.Range("M71").Value = -11266
# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
# Row 132
$ws.Range("H132").Value = 2590.6924
$ws.Range("I132").Value = 2561.7273
$ws.Range("K132").Value = 7685.1819
$ws.Range("M132").Value = -5155.1819

$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 85000
$ws.Range("J16").Value = 85000
$ws.Range("L16").Value = 85000
$ws.Range("M16").Value = -85584
# Row 81
$ws.Range("H81").Value = 1647.238
$ws.Range("I81").Value = 1341.6923
$ws.Range("K81").Value = 2683.3846
$ws.Range("M81").Value = -1622.3846
# Row 84
$ws.Range("H84").Value = 1647.238
$ws.Range("I84").Value = 1341.6923
$ws.Range("K84").Value = 13416.923
$ws.Range("M84").Value = -8112.922999999999
# Row 96
$ws.Range("H96").Value = 250874.75
$ws.Range("J96").Value = 1350
$ws.Range("L96").Value = 1350
$ws.Range("N96").Value = -4096
# Row 126
$ws.Range("H126").Value = 1535.8572
$ws.Range("I126").Value = 1423.2307
$ws.Range("K126").Value = 4269.6921
$ws.Range("M126").Value = -1799.6921
# Row 132
$ws.Range("H132").Value = 30771.334
$ws.Range("I132").Value = 2687.3872
$ws.Range("K132").Value = 8062.1616
$ws.Range("M132").Value = -5532.1616
# Row 133
$ws.Range("H133").Value = 68999
$ws.Range("J133").Value = 68999
$ws.Range("L133").Value = 68999
$ws.Range("N133").Value = -79119
